$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "test353"
$ws.Range("B2").Value = 23071746
$ws.Range("C2").Value = "narendra681"
$ws.Range("D2").Value = "w59!TP#v"
